# Apply cryptos list update (prices/volumes refreshed by GitHub Actions run)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value into a cell while forcing it to remain plain text,
# mirroring the source data (inline/shared strings), even when the text looks
# like a number (e.g. "1.001") or could be re-interpreted by Excel.
function Set-TextCell($cellRef, $text) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.ClearFormats()
}

Set-TextCell "D2" "30.312.29"
Set-TextCell "E2" "  +0.17%  "

Set-TextCell "D3" "1.868.79"
Set-TextCell "E3" "  +0.20%  "

Set-TextCell "D4" "1.001"
Set-TextCell "E4" "  +0.08%  "

Set-TextCell "D5" "234.90"
Set-TextCell "E5" "  -1.02%  "

Set-TextCell "D6" "1.001"

Set-TextCell "E7" "  +0.32%  "

Set-TextCell "E8" "  +0.08%  "

Set-TextCell "D9" "0.06574"
Set-TextCell "E9" "  +0.41%  "

Set-TextCell "D10" "21.61"
Set-TextCell "E10" "  -2.65%  "

Set-TextCell "E11" "  +1.46%  "

Set-TextCell "D12" "96.84"
Set-TextCell "E12" "  -0.97%  "

Set-TextCell "D13" "1.872.50"
Set-TextCell "E13" "  +0.30%  "

Set-TextCell "D14" "5.113"
Set-TextCell "E14" "  -1.16%  "

Set-TextCell "D15" "0.6836"
Set-TextCell "E15" "  +0.41%  "

Set-TextCell "D16" "269.81"
Set-TextCell "E16" "  -3.15%  "

Set-TextCell "D17" "30.311.24"
Set-TextCell "E17" "  +0.14%  "

Set-TextCell "D18" "14.01"
Set-TextCell "E18" "  +2.93%  "

Set-TextCell "D19" "0.000007627"
Set-TextCell "E19" "  +3.76%  "

Set-TextCell "E20" "  +0.08%  "

Set-TextCell "D21" "2.115.46"
Set-TextCell "E21" "  +0.04%  "

Set-TextCell "E22" "  +0.03%  "

Set-TextCell "D23" "5.267"
Set-TextCell "E23" "  -2.31%  "

Set-TextCell "D24" "6.203"
Set-TextCell "E24" "  +0.13%  "

Set-TextCell "D25" "9.398"
Set-TextCell "E25" "  +1.14%  "

Set-TextCell "D26" "167.99"
Set-TextCell "E26" "  -0.44%  "

Set-TextCell "D27" "18.89"
Set-TextCell "E27" "  -1.17%  "

Set-TextCell "E28" "  +0.22%  "

Set-TextCell "D29" "1.369"
Set-TextCell "E29" "  -1.21%  "

Set-TextCell "D30" "0.09883"
Set-TextCell "E30" "  +0.44%  "

Set-TextCell "D31" "4.364"
Set-TextCell "E31" "  -0.71%  "

Set-TextCell "D32" "1.463"
Set-TextCell "E32" "  -1.14%  "

Set-TextCell "D33" "4.062"
Set-TextCell "E33" "  -0.18%  "

Set-TextCell "D34" "0.04716"
Set-TextCell "E34" "  -0.64%  "

Set-TextCell "E35" "  -0.54%  "

Set-TextCell "D36" "0.6996"
Set-TextCell "E36" "  -1.26%  "

Set-TextCell "D37" "2.713"

Set-TextCell "D38" "0.01870"
Set-TextCell "E38" "  -0.43%  "

Set-TextCell "D39" "2.657"
Set-TextCell "E39" "  +1.27%  "

Set-TextCell "E40" "  -0.50%  "

Set-TextCell "D41" "71.76"
Set-TextCell "E41" "  -6.65%  "

Set-TextCell "E42" "  -0.22%  "

Set-TextCell "D43" "0.8419"
Set-TextCell "E43" "  -1.01%  "

Set-TextCell "D44" "0.4162"
Set-TextCell "E44" "  -0.65%  "

Set-TextCell "D45" "1.000"

Set-TextCell "D46" "102.69"
Set-TextCell "E46" "  -0.61%  "

Set-TextCell "D47" "7.050"
Set-TextCell "E47" "  -2.37%  "

Set-TextCell "D48" "9.153"
Set-TextCell "E48" "  -1.49%  "

Set-TextCell "D49" "910.01"
Set-TextCell "E49" "  -4.88%  "

Set-TextCell "D50" "34.42"
Set-TextCell "E50" "  +0.44%  "

Set-TextCell "D51" "0.05702"
